$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97: B97 currently holds text "3" (inline string). Convert it to a true
# numeric value of 3, leaving the rest of row 97 unchanged.
$ws.Range("B97").Value = 3

# Row 98: brand new annotation row appended after row 97.
$ws.Range("A98").Value = "Ruilin"
$ws.Range("B98").NumberFormat = "@"
$ws.Range("B98").Value = "3"
$ws.Range("B98").Style = "Normal"
$ws.Range("C98").Value = "无"
$ws.Range("D98").Value = "DIS"
$ws.Range("E98").Value = "MET"
$ws.Range("F98").Value = "2b1bc0cc-ec44-4403-95c4-ab76ff8ea3ad"
$ws.Range("G98").Value = "By3VrbbAb_annotated.xlsx"
$ws.Range("H98").Value = "The reason we need to do such sampling is because AOL dataset only consists of whole queries instead of the prefix-completion pair."
